$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the VO ID value in A13 (new term ID assigned for flu season property)
$ws.Range("A13").Value = "VO:0010415"

# Move the active selection to A14 (next empty row), as happens after entering a value
$ws.Range("A14").Select()
